$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.590.55"
$ws.Range("E2").Value = "  -3.67%  "
$ws.Range("D3").Value = "2.808.67"
$ws.Range("E3").Value = "  -3.99%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "500.66"
$ws.Range("E5").Value = "  -5.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.06"
$ws.Range("E6").Value = "  -8.55%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.523"
$ws.Range("E8").Value = "  -6.47%  "
$ws.Range("D9").Value = "2.806.42"
$ws.Range("E9").Value = "  -4.13%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.87"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.101"
$ws.Range("E11").Value = "  -6.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("E12").Value = "  -3.45%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "3.309.04"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("D15").Value = "58.790.95"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.21"
$ws.Range("E16").Value = "  -7.62%  "
$ws.Range("D17").Value = "2.816.58"
$ws.Range("E17").Value = "  -3.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000131"
$ws.Range("E18").Value = "  -6.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.64"
$ws.Range("E19").Value = "  -7.07%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.68"
$ws.Range("E20").Value = "  -5.19%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.83"
$ws.Range("E21").Value = "  -7.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.63"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.57"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.421"
$ws.Range("E26").Value = "  -7.12%  "
$ws.Range("E27").Value = "  -8.33%  "
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  -6.96%  "
$ws.Range("D30").Value = "0.0₃0780"
$ws.Range("E30").Value = "  -10.59%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.59"
$ws.Range("E32").Value = "  -5.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.80"
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.85"
$ws.Range("E34").Value = "  -3.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.07"
$ws.Range("E35").Value = "  -7.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.22"
$ws.Range("E36").Value = "  -7.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.886"
$ws.Range("E37").Value = "  -12.27%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.85"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.10"
$ws.Range("E39").Value = "  -9.47%  "
$ws.Range("D40").Value = "2.208.01"
$ws.Range("E40").Value = "  -6.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.622"
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.46"
$ws.Range("E43").Value = "  -6.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0549"
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.32"
$ws.Range("E45").Value = "  -11.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.85"
$ws.Range("E46").Value = "  -10.65%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0222"
$ws.Range("E48").Value = "  -5.47%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0882"
$ws.Range("E49").Value = "  -5.10%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.51"
$ws.Range("E50").Value = "  -8.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.06"
$ws.Range("E51").Value = "  -8.76%  "
